# Auto-generated edit script: updates market price / profit data cells
# as scraped by the scheduled runner, across multiple Leve-profit sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 25643832
$ws.Range("I15").Value = 25643832
$ws.Range("K15").Value = 76931496
$ws.Range("M15").Value = -76931327

# Row 32
$ws.Range("H32").Value = 433.33334
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 433.33334
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 433.33334
$ws.Range("M32").ClearContents() | Out-Null
$ws.Range("N32").Value = -1085.33334

# Row 100
$ws.Range("H100").Value = 5945.75
$ws.Range("I100").Value = 5698.4287
$ws.Range("J100").Value = 6292
$ws.Range("K100").Value = 5698.4287
$ws.Range("L100").Value = 6292
$ws.Range("M100").Value = -5157.4287
$ws.Range("N100").Value = -7374

# Row 113
$ws.Range("H113").Value = 4436.95
$ws.Range("I113").Value = 4299.8887
$ws.Range("J113").Value = 4549.091
$ws.Range("K113").Value = 4299.8887
$ws.Range("L113").Value = 4549.091
$ws.Range("M113").Value = -1045.8887
$ws.Range("N113").Value = -11057.091

# Row 132
$ws.Range("H132").Value = 4350208
$ws.Range("I132").Value = 4653423
$ws.Range("J132").Value = 4130
$ws.Range("K132").Value = 13960269
$ws.Range("L132").Value = 12390
$ws.Range("M132").Value = -13957739
$ws.Range("N132").Value = -17450

# Row 137
$ws.Range("H137").Value = 3033699.5
$ws.Range("I137").Value = 3849484
$ws.Range("K137").Value = 11548452
$ws.Range("M137").Value = -11545902

# Row 141
$ws.Range("H141").Value = 726964.1
$ws.Range("I141").Value = 2195
$ws.Range("J141").Value = 2418092.2
$ws.Range("K141").Value = 6585
$ws.Range("L141").Value = 7254276.600000001
$ws.Range("M141").Value = -1405
$ws.Range("N141").Value = -7264636.600000001


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 88
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents() | Out-Null
$ws.Range("N88").ClearContents() | Out-Null

# Row 91
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents() | Out-Null
$ws.Range("N91").ClearContents() | Out-Null

# Row 122
$ws.Range("H122").Value = 5498.5713
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 5998.3335
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 17995.0005
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -22895.0005


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 2669.0908
$ws.Range("I99").Value = 1742.7142
$ws.Range("J99").Value = 4290.25
$ws.Range("K99").Value = 1742.7142
$ws.Range("L99").Value = 4290.25
$ws.Range("M99").Value = -244.7141999999999
$ws.Range("N99").Value = -7286.25

# Row 105
$ws.Range("H105").Value = 2011.1428
$ws.Range("I105").Value = 1869.1666
$ws.Range("J105").Value = 2200.4443
$ws.Range("K105").Value = 1869.1666
$ws.Range("L105").Value = 2200.4443
$ws.Range("M105").Value = -122.1666
$ws.Range("N105").Value = -5694.4443

# Row 107
$ws.Range("H107").Value = 3232.5
$ws.Range("I107").Value = 2227.7144
$ws.Range("J107").Value = 4014
$ws.Range("K107").Value = 2227.7144
$ws.Range("L107").Value = 4014
$ws.Range("M107").Value = -307.7143999999998
$ws.Range("N107").Value = -7854


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2085799.8
$ws.Range("J31").Value = 4565
$ws.Range("L31").Value = 4565
$ws.Range("N31").Value = -5155

# Row 34
$ws.Range("H34").Value = 2085799.8
$ws.Range("J34").Value = 4565
$ws.Range("L34").Value = 4565
$ws.Range("N34").Value = -4969

# Row 86
$ws.Range("H86").Value = 3439.8108
$ws.Range("I86").Value = 2717.3809
$ws.Range("J86").Value = 4388
$ws.Range("K86").Value = 2717.3809
$ws.Range("L86").Value = 4388
$ws.Range("M86").Value = -1594.3809
$ws.Range("N86").Value = -6634

# Row 89
$ws.Range("H89").Value = 3439.8108
$ws.Range("I89").Value = 2717.3809
$ws.Range("J89").Value = 4388
$ws.Range("K89").Value = 13586.9045
$ws.Range("L89").Value = 21940
$ws.Range("M89").Value = -7970.904500000001
$ws.Range("N89").Value = -33172

# Row 99
$ws.Range("H99").Value = 3800
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 3800
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 3800
$ws.Range("M99").ClearContents() | Out-Null
$ws.Range("N99").Value = -6796

# Row 126
$ws.Range("H126").Value = 3800
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 3800
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 11400
$ws.Range("M126").ClearContents() | Out-Null
$ws.Range("N126").Value = -16340

# Row 134
$ws.Range("H134").Value = 2079.7896
$ws.Range("I134").Value = 1076.3334
$ws.Range("J134").Value = 3800
$ws.Range("K134").Value = 3229.0002
$ws.Range("L134").Value = 11400
$ws.Range("M134").Value = -694.0001999999999
$ws.Range("N134").Value = -16470


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 3201.6667
$ws.Range("I97").Value = 1852.5
$ws.Range("J97").Value = 5900
$ws.Range("K97").Value = 1852.5
$ws.Range("L97").Value = 5900
$ws.Range("M97").Value = -1356.5
$ws.Range("N97").Value = -6892

# Row 122
$ws.Range("H122").Value = 3588.8
$ws.Range("I122").Value = 2463
$ws.Range("J122").Value = 3998.182
$ws.Range("K122").Value = 7389
$ws.Range("L122").Value = 11994.546
$ws.Range("M122").Value = -4939
$ws.Range("N122").Value = -16894.546

# Row 132
$ws.Range("H132").Value = 3329.7856
$ws.Range("I132").Value = 3071.3044
$ws.Range("J132").Value = 3642.6843
$ws.Range("K132").Value = 9213.913199999999
$ws.Range("L132").Value = 10928.0529
$ws.Range("M132").Value = -6683.913199999999
$ws.Range("N132").Value = -15988.0529


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 1596.6666
$ws.Range("I40").Value = 1393.3334
$ws.Range("K40").Value = 1393.3334
$ws.Range("M40").Value = -1257.3334

# Row 60
$ws.Range("H60").Value = 14000
$ws.Range("J60").Value = 14000
$ws.Range("L60").Value = 14000
$ws.Range("N60").Value = -15018

# Row 93
$ws.Range("H93").Value = 3012.5625
$ws.Range("I93").Value = 1771.5714
$ws.Range("K93").Value = 1771.5714
$ws.Range("M93").Value = -523.5714

# Row 100
$ws.Range("H100").Value = 2225.8572
$ws.Range("I100").Value = 1650
$ws.Range("J100").Value = 2382.9092
$ws.Range("K100").Value = 1650
$ws.Range("L100").Value = 2382.9092
$ws.Range("M100").Value = -1109
$ws.Range("N100").Value = -3464.9092

# Row 132
$ws.Range("H132").Value = 3088.9033
$ws.Range("I132").Value = 1722.0625
$ws.Range("J132").Value = 4546.8667
$ws.Range("K132").Value = 5166.1875
$ws.Range("L132").Value = 13640.6001
$ws.Range("M132").Value = -2636.1875
$ws.Range("N132").Value = -18700.6001

# Row 136
$ws.Range("H136").Value = 4549013
$ws.Range("I136").Value = 10003970
$ws.Range("J136").Value = 3214.9167
$ws.Range("K136").Value = 30011910
$ws.Range("L136").Value = 9644.750100000001
$ws.Range("M136").Value = -30009360
$ws.Range("N136").Value = -14744.7501


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 627847.0600000001
$ws.Range("I122").Value = 771504.1
$ws.Range("K122").Value = 2314512.3
$ws.Range("M122").Value = -2312062.3

# Row 132
$ws.Range("H132").Value = 226625.69
$ws.Range("I132").Value = 358836
$ws.Range("K132").Value = 1076508
$ws.Range("M132").Value = -1073978

# Row 136
$ws.Range("H136").Value = 879.5
$ws.Range("I136").Value = 537.72095
$ws.Range("J136").Value = 2010
$ws.Range("K136").Value = 1613.16285
$ws.Range("L136").Value = 6030
$ws.Range("M136").Value = 936.8371499999998
$ws.Range("N136").Value = -11130

